$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update O4:O9 from 2 to 0.5 (calibration constant changed)
$ws.Range("O4").Value = 0.5
$ws.Range("O5").Value = 0.5
$ws.Range("O6").Value = 0.5
$ws.Range("O7").Value = 0.5
$ws.Range("O8").Value = 0.5
$ws.Range("O9").Value = 0.5

# 2. Update row 35 header labels from temperature/voltage units to percentage units
$ws.Range("H35").Value = "Sx (%)"

$ws.Range("I35").Value = "Sx2 (% 2)"
$ws.Range("I35").Characters(3,1).Font.Superscript = $true
$ws.Range("I35").Characters(8,1).Font.Superscript = $true

$ws.Range("J35").Value = "Sy (%)"
$ws.Range("K35").Value = "Sxy (% . %)"

# 3. Delete the now-redundant duplicate MMQ block (rows 61-68)
$ws.Range("G61:K68").Select()
$ws.Range("A61:A68").EntireRow.Delete()
